$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.135.17'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.559.66'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.29%  '
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.99'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.63%  '
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0595'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0870'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.783.13'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.582.28'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.76'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.516'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.136.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.74'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '216.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("E19").Value = '  +0.93%  '
$ws.Range("E20").Value = '  -0.73%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.11'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.19'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("E24").Value = '  +0.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.98'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.70%  '
$ws.Range("E28").Value = '  +1.41%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("E30").Value = '  +2.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0469'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.24%  '
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("E33").Value = '  +1.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.433.44'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.84%  '
$ws.Range("E35").Value = '  +4.04%  '
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("E37").Value = '  -0.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0165'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.35%  '
$ws.Range("E39").Value = '  -0.60%  '
$ws.Range("E40").Value = '  +1.54%  '
$ws.Range("E41").Value = '  -0.33%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E43").Value = '  +0.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.997'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.21'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.63%  '
$ws.Range("E46").Value = '  -1.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.695.63'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0524'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₇0999'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0948'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.09%  '

Write-Output "Updated $($ws.Name) with latest crypto prices"
